$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue 'D2' '37.423.97'
Set-TextValue 'E2' '  -1.08%  '
Set-TextValue 'D3' '2.051.32'
Set-TextValue 'E3' '  -1.81%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '228.93'
Set-TextValue 'E5' '  -2.09%  '
Set-TextValue 'D6' '0.611'
Set-TextValue 'E6' '  -2.18%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '56.31'
Set-TextValue 'E8' '  -3.65%  '
Set-TextValue 'D9' '0.386'
Set-TextValue 'E9' '  -2.22%  '
Set-TextValue 'D10' '0.0809'
Set-TextValue 'E10' '  +3.35%  '
Set-TextValue 'E11' '  -2.08%  '
Set-TextValue 'D12' '2.352.25'
Set-TextValue 'E12' '  -1.89%  '
Set-TextValue 'D13' '14.54'
Set-TextValue 'E13' '  -4.96%  '
Set-TextValue 'D14' '20.63'
Set-TextValue 'E14' '  -2.90%  '
Set-TextValue 'D15' '0.754'
Set-TextValue 'E15' '  -3.32%  '
Set-TextValue 'D16' '5.26'
Set-TextValue 'E16' '  -2.11%  '
Set-TextValue 'D17' '2.039.45'
Set-TextValue 'E17' '  -2.39%  '
Set-TextValue 'D18' '37.309.95'
Set-TextValue 'E18' '  -1.37%  '
Set-TextValue 'D19' '6.11'
Set-TextValue 'E19' '  -0.68%  '
Set-TextValue 'D20' '69.82'
Set-TextValue 'E20' '  -1.76%  '
Set-TextValue 'D21' '0.0₃0849'
Set-TextValue 'E21' '  +1.47%  '
Set-TextValue 'D22' '225.81'
Set-TextValue 'E22' '  -1.80%  '
Set-TextValue 'D23' '1.00'
Set-TextValue 'E23' '  +0.31%  '
Set-TextValue 'E24' '  -1.01%  '
Set-TextValue 'D25' '2.28'
Set-TextValue 'E25' '  -4.65%  '
Set-TextValue 'D26' '9.55'
Set-TextValue 'E26' '  -2.38%  '
Set-TextValue 'D27' '168.62'
Set-TextValue 'E27' '  -1.77%  '
Set-TextValue 'E28' '  -3.73%  '
Set-TextValue 'D29' '1.38'
Set-TextValue 'E29' '  -1.80%  '
Set-TextValue 'D30' '18.90'
Set-TextValue 'E30' '  -3.14%  '
Set-TextValue 'E31' '  -2.57%  '
Set-TextValue 'D32' '4.55'
Set-TextValue 'E32' '  -3.27%  '
Set-TextValue 'D33' '0.0613'
Set-TextValue 'E33' '  -3.20%  '
Set-TextValue 'D34' '4.54'
Set-TextValue 'E34' '  -2.10%  '
Set-TextValue 'D35' '2.40'
Set-TextValue 'E35' '  -3.85%  '
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'E37' '  -0.01%  '
Set-TextValue 'D38' '3.19'
Set-TextValue 'E38' '  -4.04%  '
Set-TextValue 'D39' '5.43'
Set-TextValue 'E39' '  +0.52%  '
Set-TextValue 'B40' 'VeChain'
Set-TextValue 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D40' '0.0221'
Set-TextValue 'E40' '  -6.15%  '
Set-TextValue 'B41' 'Maker'
Set-TextValue 'C41' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D41' '1.502.20'
Set-TextValue 'E41' '  +3.34%  '
Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '16.90'
Set-TextValue 'E42' '  +0.60%  '
Set-TextValue 'B43' 'HuobiToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D43' '2.88'
Set-TextValue 'E43' '  -1.60%  '
Set-TextValue 'D44' '96.15'
Set-TextValue 'E44' '  -5.04%  '
Set-TextValue 'D45' '0.0934'
Set-TextValue 'E45' '  -3.99%  '
Set-TextValue 'D46' '1.14'
Set-TextValue 'E46' '  -4.20%  '
Set-TextValue 'D47' '1.02'
Set-TextValue 'E47' '  -4.31%  '
Set-TextValue 'D48' '7.21'
Set-TextValue 'E48' '  -0.19%  '
Set-TextValue 'B49' 'FTXToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D49' '3.85'
Set-TextValue 'E49' '  -6.55%  '
Set-TextValue 'B50' 'MXToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D50' '2.92'
Set-TextValue 'E50' '  -1.64%  '
Set-TextValue 'D51' '2.237.19'
Set-TextValue 'E51' '  -1.95%  '
